$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.158182694426159
$ws.Range("C2").Value = 0.339459061613951
$ws.Range("D2").Value = 0.2219658837373615
$ws.Range("F2").Value = 1.255498208774078
$ws.Range("G2").Value = 0.6524449674027863
$ws.Range("H2").Value = 0.7760300657686301
$ws.Range("J2").Value = 0.2335198683620625
$ws.Range("L2").Value = 0.3310972822870113
$ws.Range("M2").Value = 0.2914151296881684
$ws.Range("O2").Value = 2.848292395402837
$ws.Range("B3").Value = 1.051963321756546
$ws.Range("C3").Value = 0.3291491866946217
$ws.Range("D3").Value = 0.2220636002734437
$ws.Range("F3").Value = 1.266306239975094
$ws.Range("G3").Value = 0.6589094181134953
$ws.Range("H3").Value = 0.7834690126521906
$ws.Range("J3").Value = 0.2357153974068122
$ws.Range("L3").Value = 0.3270628312532722
$ws.Range("M3").Value = 0.27392787821897
$ws.Range("O3").Value = 2.877285259570826
$ws.Range("B4").Value = 0.9866901134113846
$ws.Range("C4").Value = 0.3228437193578202
$ws.Range("D4").Value = 0.2221918684999054
$ws.Range("F4").Value = 1.273683932259544
$ws.Range("G4").Value = 0.6633787005981162
$ws.Range("H4").Value = 0.7884179818569805
$ws.Range("J4").Value = 0.2371419203373524
$ws.Range("L4").Value = 0.324680761034557
$ws.Range("M4").Value = 0.2632220308190938
$ws.Range("O4").Value = 2.896935745553932
$ws.Range("B5").Value = 0.9600790786992093
$ws.Range("C5").Value = 0.320280701021872
$ws.Range("D5").Value = 0.2222613602220136
$ws.Range("F5").Value = 1.276876875030744
$ws.Range("G5").Value = 0.6653255951242869
$ws.Range("H5").Value = 0.7905306749387364
$ws.Range("J5").Value = 0.2377429962928801
$ws.Range("L5").Value = 0.3237340992758391
$ws.Range("M5").Value = 0.258867545416912
$ws.Range("O5").Value = 2.905408164926527
$ws.Range("B6").Value = 0.9556596886456532
$ws.Range("C6").Value = 0.3198555148613593
$ws.Range("D6").Value = 0.2222739409972547
$ws.Range("F6").Value = 1.277418323634912
$ws.Range("G6").Value = 0.6656564592904601
$ws.Range("H6").Value = 0.7908872820732284
$ws.Range("J6").Value = 0.2378439985270582
$ws.Range("L6").Value = 0.3235783635068259
$ws.Range("M6").Value = 0.2581449937991067
$ws.Range("O6").Value = 2.906843063854964
$ws.Range("B7").Value = 0.9863312723885542
$ws.Range("C7").Value = 0.3228091268817792
$ws.Range("D7").Value = 0.2221927358875817
$ws.Range("F7").Value = 1.273726238363167
$ws.Range("G7").Value = 0.6634044486115513
$ws.Range("H7").Value = 0.788446085839162
$ws.Range("J7").Value = 0.2371499466295841
$ws.Range("L7").Value = 0.3246678964568943
$ws.Range("M7").Value = 0.2631632709138856
$ws.Range("O7").Value = 2.897048126328997
$ws.Range("B8").Value = 1.121570599767608
$ws.Range("C8").Value = 0.3358992166601524
$ws.Range("D8").Value = 0.2219854392261738
$ws.Range("F8").Value = 1.259070923315313
$ws.Range("G8").Value = 0.6545700427549122
$ws.Range("H8").Value = 0.7785158712093505
$ws.Range("J8").Value = 0.2342606189641447
$ws.Range("L8").Value = 0.3296865465042771
$ws.Range("M8").Value = 0.2853792067719567
$ws.Range("O8").Value = 2.857905309046501
$ws.Range("B9").Value = 1.386273960168012
$ws.Range("C9").Value = 0.3617556779333029
$ws.Range("D9").Value = 0.2221183787989958
$ws.Range("F9").Value = 1.236215368041869
$ws.Range("G9").Value = 0.6412196067032809
$ws.Range("H9").Value = 0.7620676390600494
$ws.Range("J9").Value = 0.2292159318383931
$ws.Range("L9").Value = 0.3402773441994498
$ws.Range("M9").Value = 0.3291810491859337
$ws.Range("O9").Value = 2.795824855126057
$ws.Range("B10").Value = 1.580369910084755
$ws.Range("C10").Value = 0.3808546300785736
$ws.Range("D10").Value = 0.2225421095002105
$ws.Range("F10").Value = 1.223010564959736
$ws.Range("G10").Value = 0.6338420701466134
$ws.Range("H10").Value = 0.7518250896088148
$ws.Range("J10").Value = 0.2258865391454581
$ws.Range("L10").Value = 0.3485090662048407
$ws.Range("M10").Value = 0.361492752882846
$ws.Range("O10").Value = 2.759176799068086
$ws.Range("B11").Value = 1.668571470719542
$ws.Range("C11").Value = 0.389563118399991
$ws.Range("D11").Value = 0.2228050979172664
$ws.Range("F11").Value = 1.217782421092792
$ws.Range("G11").Value = 0.6310155496214946
$ws.Range("H11").Value = 0.7475650238031335
$ws.Range("J11").Value = 0.2244533778432984
$ws.Range("L11").Value = 0.3523505282265802
$ws.Range("M11").Value = 0.376217852059419
$ws.Range("O11").Value = 2.744453844133204
$ws.Range("B12").Value = 1.701955981466256
$ws.Range("C12").Value = 0.3928634797169934
$ws.Range("D12").Value = 0.2229147322895955
$ws.Range("F12").Value = 1.215914654314631
$ws.Range("G12").Value = 0.630021513030897
$ws.Range("H12").Value = 0.7460092403944429
$ws.Range("J12").Value = 0.2239223546319664
$ws.Range("L12").Value = 0.353818987807287
$ws.Range("M12").Value = 0.3817973581640786
$ws.Range("O12").Value = 2.739159075542602
$ws.Range("B13").Value = 1.694766751690736
$ws.Range("C13").Value = 0.3921525744593737
$ws.Range("D13").Value = 0.2228906745906656
$ws.Range("F13").Value = 1.216311928835715
$ws.Range("G13").Value = 0.6302322004337668
$ws.Range("H13").Value = 0.7463417529903325
$ws.Range("J13").Value = 0.2240362006661893
$ws.Range("L13").Value = 0.3535021180586995
$ws.Range("M13").Value = 0.3805955647862689
$ws.Range("O13").Value = 2.74028691866306
$ws.Range("B14").Value = 1.671318356898894
$ws.Range("C14").Value = 0.3898345896240869
$ws.Range("D14").Value = 0.2228139165612291
$ws.Range("F14").Value = 1.217626513629213
$ws.Range("G14").Value = 0.630932239473907
$ws.Range("H14").Value = 0.7474358780256267
$ws.Range("J14").Value = 0.2244094561999379
$ws.Range("L14").Value = 0.3524710638819641
$ws.Range("M14").Value = 0.3766768146700699
$ws.Range("O14").Value = 2.744012617018967
$ws.Range("B15").Value = 1.656953458436533
$ws.Range("C15").Value = 0.3884150935668345
$ws.Range("D15").Value = 0.2227682068204402
$ws.Range("F15").Value = 1.218446323644962
$ws.Range("G15").Value = 0.6313709754027172
$ws.Range("H15").Value = 0.7481135374920456
$ws.Range("J15").Value = 0.224639606964498
$ws.Range("L15").Value = 0.3518413040142292
$ws.Range("M15").Value = 0.3742769032667113
$ws.Range("O15").Value = 2.746331252784131
$ws.Range("B16").Value = 1.574603671474506
$ws.Range("C16").Value = 0.3802858947844072
$ws.Range("D16").Value = 0.2225263313038965
$ws.Range("F16").Value = 1.22336791017306
$ws.Range("G16").Value = 0.634037459309198
$ws.Range("H16").Value = 0.7521115293675393
$ws.Range("J16").Value = 0.22598183572015
$ws.Range("L16").Value = 0.34825995372573
$ws.Range("M16").Value = 0.3605309325053554
$ws.Range("O16").Value = 2.760178213035772
$ws.Range("B17").Value = 1.524059301582213
$ws.Range("C17").Value = 0.3753038953992416
$ws.Range("D17").Value = 0.2223958987960017
$ws.Range("F17").Value = 1.226586631394341
$ws.Range("G17").Value = 0.6358090045470988
$ws.Range("H17").Value = 0.7546664359591659
$ws.Range("J17").Value = 0.2268260821925292
$ws.Range("L17").Value = 0.3460876070470391
$ws.Range("M17").Value = 0.3521047167528479
$ws.Range("O17").Value = 2.769172138208035
$ws.Range("B18").Value = 1.494978822044857
$ws.Range("C18").Value = 0.3724403058801613
$ws.Range("D18").Value = 0.2223274911091977
$ws.Range("F18").Value = 1.228511256194594
$ws.Range("G18").Value = 0.6368777858072718
$ws.Range("H18").Value = 0.7561735373211462
$ws.Range("J18").Value = 0.2273193330454415
$ws.Range("L18").Value = 0.3448472512776846
$ws.Range("M18").Value = 0.3472606810894163
$ws.Range("O18").Value = 2.774528583017855
$ws.Range("B19").Value = 1.485131240314388
$ws.Range("C19").Value = 0.3714710825329917
$ws.Range("D19").Value = 0.2223054669512621
$ws.Range("F19").Value = 1.22917548977378
$ws.Range("G19").Value = 0.6372482112699629
$ws.Range("H19").Value = 0.75669027182375
$ws.Range("J19").Value = 0.227487656008778
$ws.Range("L19").Value = 0.344428858868028
$ws.Range("M19").Value = 0.3456210142830045
$ws.Range("O19").Value = 2.776373668760939
$ws.Range("B20").Value = 1.529440750729009
$ws.Range("C20").Value = 0.3758340405010756
$ws.Range("D20").Value = 0.2224090994553762
$ws.Range("F20").Value = 1.226236406380693
$ws.Range("G20").Value = 0.6356152616331912
$ws.Range("H20").Value = 0.7543905717895143
$ws.Range("J20").Value = 0.2267354178142202
$ws.Range("L20").Value = 0.3463179141289601
$ws.Range("M20").Value = 0.3530014455801123
$ws.Range("O20").Value = 2.768195738412714
$ws.Range("B21").Value = 1.678206155849352
$ws.Range("C21").Value = 0.3905153682950413
$ws.Range("D21").Value = 0.2228361900197626
$ws.Range("F21").Value = 1.217237347851693
$ws.Range("G21").Value = 0.6307245489005737
$ws.Range("H21").Value = 0.7471129489480575
$ws.Range("J21").Value = 0.2242995050687497
$ws.Range("L21").Value = 0.3527735365790647
$ws.Range("M21").Value = 0.377827756303617
$ws.Range("O21").Value = 2.742910674071936
$ws.Range("B22").Value = 1.775341630647233
$ws.Range("C22").Value = 0.4001257628733583
$ws.Range("D22").Value = 0.2231738518958082
$ws.Range("F22").Value = 1.212008881156024
$ws.Range("G22").Value = 0.627973022630556
$ws.Range("H22").Value = 0.7426912321786148
$ws.Range("J22").Value = 0.2227755883878295
$ws.Range("L22").Value = 0.3570729223556981
$ws.Range("M22").Value = 0.3940730173904328
$ws.Range("O22").Value = 2.72802050445793
$ws.Range("B23").Value = 1.723507625034927
$ws.Range("C23").Value = 0.3949952032470776
$ws.Range("D23").Value = 0.2229882969876087
$ws.Range("F23").Value = 1.214739659047268
$ws.Range("G23").Value = 0.6294008103695461
$ws.Range("H23").Value = 0.745020570540774
$ws.Range("J23").Value = 0.2235827076644412
$ws.Range("L23").Value = 0.3547709612455634
$ws.Range("M23").Value = 0.3854009187726035
$ws.Range("O23").Value = 2.735817955674719
$ws.Range("B24").Value = 1.527007865470011
$ws.Range("C24").Value = 0.375594359923582
$ws.Range("D24").Value = 0.2224031109389699
$ws.Range("F24").Value = 1.226394512173705
$ws.Range("G24").Value = 0.6357026961471206
$ws.Range("H24").Value = 0.754515170828185
$ws.Range("J24").Value = 0.2267763826249567
$ws.Range("L24").Value = 0.3462137656386375
$ws.Range("M24").Value = 0.3525960334910749
$ws.Range("O24").Value = 2.76863659028561
$ws.Range("B25").Value = 1.314726232788644
$ws.Range("C25").Value = 0.3547420373204204
$ws.Range("D25").Value = 0.2220249536966818
$ws.Range("F25").Value = 1.241768484169654
$ws.Range("G25").Value = 0.6444049468326227
$ws.Range("H25").Value = 0.7661937021451166
$ws.Range("J25").Value = 0.2305143296117862
$ws.Range("L25").Value = 0.3373326736925151
$ws.Range("M25").Value = 0.3173075963218679
$ws.Range("O25").Value = 2.811046364548801
